$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 19, shifting existing rows 19:25 down to 20:26
$ws.Rows.Item(19).Insert()

# Populate the newly inserted row 19 with this week's data
$ws.Cells.Item(19, 1).Value = 10
$ws.Cells.Item(19, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(19, 3).Value = "La Araucanía"
$ws.Cells.Item(19, 4).Value = 44455
$ws.Cells.Item(19, 4).Style = $ws.Cells.Item(20, 4).Style
$ws.Cells.Item(19, 4).NumberFormat = $ws.Cells.Item(20, 4).NumberFormat
$ws.Cells.Item(19, 5).Value = 9
$ws.Cells.Item(19, 6).Value = 100114002
$ws.Cells.Item(19, 7).Value = "Camote"
$ws.Cells.Item(19, 8).Value = "Sin especificar"
$ws.Cells.Item(19, 9).Value = "Primera"
$ws.Cells.Item(19, 10).Value = 30
$ws.Cells.Item(19, 11).Value = 20000
$ws.Cells.Item(19, 12).Value = 20000
$ws.Cells.Item(19, 13).Value = 20000
$ws.Cells.Item(19, 14).Value = "$/malla 20 kilos"
$ws.Cells.Item(19, 15).Value = "Perú"
$ws.Cells.Item(19, 16).Value = 1000
$ws.Cells.Item(19, 17).Value = 20
$ws.Cells.Item(19, 18).Value = "Hortaliza"
